$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings)
$ws.Range("A1").Value = "tiempo"
$ws.Range("B1").Value = "Temperatura"

# Update selected cell
$ws.Range("B7").Select()
